$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$HAlignLeft   = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignLeft
$HAlignCenter = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignCenter
$VAlignCenter = [Microsoft.Office.Interop.Excel.XlVAlign]::xlVAlignCenter
$PasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# ---------------------------------------------------------------------
# 1) Preserve the formatting currently on row 42 (the "highlighted" /
#    final row of the log) so it can be moved down to the new last row
#    (44) before we overwrite row 42's own formatting.
# ---------------------------------------------------------------------
$ws.Range("A42").Copy()
$ws.Range("A44").PasteSpecial($PasteFormats)
$ws.Range("B42").Copy()
$ws.Range("B44").PasteSpecial($PasteFormats)
$ws.Range("C42").Copy()
$ws.Range("C44").PasteSpecial($PasteFormats)

# ---------------------------------------------------------------------
# 2) Row 42 becomes a normal (non-highlighted) data row, matching the
#    date/hours formatting already used elsewhere in the log.
# ---------------------------------------------------------------------
$ws.Range("A41").Copy()
$ws.Range("A42").PasteSpecial($PasteFormats)
$ws.Range("B41").Copy()
$ws.Range("B42").PasteSpecial($PasteFormats)

# ---------------------------------------------------------------------
# 3) New row 43 (merged together with row 42 in column C) gets the same
#    plain-row formatting.
# ---------------------------------------------------------------------
$ws.Range("A41").Copy()
$ws.Range("A43").PasteSpecial($PasteFormats)
$ws.Range("B41").Copy()
$ws.Range("B43").PasteSpecial($PasteFormats)

# ---------------------------------------------------------------------
# 3b) Column C of rows 42-43 is merged into a single note cell. Merge
#     first, then format the merged range as a single unit - copying
#     format piecemeal into cells that are already split across a
#     merge boundary produces a partial (top/bottom-less) synthetic
#     border instead of the plain thin border used everywhere else.
# ---------------------------------------------------------------------
$ws.Range("C42:C43").Merge()
$ws.Range("B41").Copy()
$ws.Range("C42:C43").PasteSpecial($PasteFormats)
$ws.Range("C42:C43").HorizontalAlignment = $HAlignLeft
$ws.Range("C42:C43").VerticalAlignment = $VAlignCenter
$ws.Range("C42:C43").WrapText = $true

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 4) Write the cell values/dates/hours/notes for rows 42-44.
# ---------------------------------------------------------------------
$ws.Range("A42").Value2 = 45744
$ws.Range("B42").Value2 = 5
$ws.Range("C42").Value = "With simulated environment, test results were obtained and analyzed for validation phase"

$ws.Range("A43").Value2 = 45745
$ws.Range("B43").Value2 = 5

$ws.Range("A44").Value2 = 45746
$ws.Range("B44").Value2 = 4
$ws.Range("C44").Value = "Compared known vulnerabilities from OWASP Juice Shop with test results, started final report generation"

# ---------------------------------------------------------------------
# 6) Rows 37-40 are now collapsed (hidden) in the log view.
# ---------------------------------------------------------------------
$ws.Rows.Item(37).Hidden = $true
$ws.Rows.Item(38).Hidden = $true
$ws.Rows.Item(39).Hidden = $true
$ws.Rows.Item(40).Hidden = $true

# ---------------------------------------------------------------------
# 7) Update the view: scroll back to the top-left and move the
#    selection down to the new empty row under the log.
# ---------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B50").Select()
